# Scheduled-runner refresh of the FFXIV leve-profit market data.
# Columns H-N on each sheet are plain cached numbers (no formulas),
# pulled from an external price feed, so the edit is a straight
# per-cell overwrite; a few rows also gain or lose a trailing
# cell (M/N) because the refreshed row is one field wider/narrower.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 6636.6665
$ws.Range("I76").Value = 6514.636
$ws.Range("J76").Value = 6828.4287
$ws.Range("K76").Value = 6514.636
$ws.Range("L76").Value = 6828.4287
$ws.Range("M76").Value = -6199.636
$ws.Range("N76").Value = -7458.4287
# Row 79
$ws.Range("H79").Value = 6636.6665
$ws.Range("I79").Value = 6514.636
$ws.Range("J79").Value = 6828.4287
$ws.Range("K79").Value = 6514.636
$ws.Range("L79").Value = 6828.4287
$ws.Range("M79").Value = -5422.636
$ws.Range("N79").Value = -9012.4287
# Row 116
$ws.Range("H116").Value = 9854.909
$ws.Range("I116").Value = 11050.5
$ws.Range("J116").Value = 6666.6665
$ws.Range("K116").Value = 11050.5
$ws.Range("L116").Value = 6666.6665
$ws.Range("M116").Value = -7608.5
$ws.Range("N116").Value = -13550.6665
# Row 138
$ws.Range("H138").Value = 2790.7805
$ws.Range("J138").Value = 3519.5789
$ws.Range("L138").Value = 10558.7367
$ws.Range("N138").Value = -20838.7367

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 25000950
$ws.Range("J102").Value = 1450
$ws.Range("L102").Value = 1450
$ws.Range("N102").Value = -4694

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3718.15
$ws.Range("I86").Value = 3412.0715
$ws.Range("J86").Value = 4432.3335
$ws.Range("K86").Value = 3412.0715
$ws.Range("L86").Value = 4432.3335
$ws.Range("M86").Value = -2289.0715
$ws.Range("N86").Value = -6678.3335
# Row 89
$ws.Range("H89").Value = 3718.15
$ws.Range("I89").Value = 3412.0715
$ws.Range("J89").Value = 4432.3335
$ws.Range("K89").Value = 17060.3575
$ws.Range("L89").Value = 22161.6675
$ws.Range("M89").Value = -11444.3575
$ws.Range("N89").Value = -33393.6675

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10998.223
$ws.Range("I31").Value = 7716.9473
$ws.Range("K31").Value = 7716.9473
$ws.Range("M31").Value = -7421.9473
# Row 34
$ws.Range("H34").Value = 10998.223
$ws.Range("I34").Value = 7716.9473
$ws.Range("K34").Value = 7716.9473
$ws.Range("M34").Value = -7514.9473
# Row 69
$ws.Range("H69").Value = 44999
$ws.Range("I69").Value = 44999
$ws.Range("K69").Value = 44999
$ws.Range("M69").Value = -44250
# Row 72
$ws.Range("H72").Value = 44999
$ws.Range("I72").Value = 44999
$ws.Range("K72").Value = 134997
$ws.Range("M72").Value = -131253
# Row 125
$ws.Range("H125").Value = 54999
$ws.Range("J125").Value = 54999
$ws.Range("L125").Value = 54999
$ws.Range("N125").Value = -59919

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 115263.34
$ws.Range("I11").Value = 115263.34
$ws.Range("K11").Value = 345790.02
$ws.Range("M11").Value = -345650.02
# Row 23
$ws.Range("H23").Value = 1999.5
$ws.Range("I23").Value = 1999
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 5997
$ws.Range("L23").Value = 6000
$ws.Range("M23").Value = -5762
$ws.Range("N23").Value = -6470
# Row 92
$ws.Range("H92").Value = 429.4
$ws.Range("I92").Value = 436.75
$ws.Range("K92").Value = 1310.25
$ws.Range("M92").Value = -62.25
# Row 107
$ws.Range("H107").Value = 1903
$ws.Range("I107").Value = 592.2857
$ws.Range("J107").Value = 3049.875
$ws.Range("K107").Value = 1776.8571
$ws.Range("L107").Value = 9149.625
$ws.Range("M107").Value = 143.1428999999998
$ws.Range("N107").Value = -12989.625
# Row 109
$ws.Range("H109").Value = 2013.5
$ws.Range("I109").Value = 2013.5
$ws.Range("K109").Value = 6040.5
$ws.Range("M109").Value = -5000.5
# Row 110
$ws.Range("H110").Value = 4500
$ws.Range("I110").Value = 4500
$ws.Range("K110").Value = 13500
$ws.Range("M110").Value = -9410
# Row 131
$ws.Range("H131").Value = 1774.25
$ws.Range("I131").Value = 1580
$ws.Range("J131").Value = 1913
$ws.Range("K131").Value = 4740
$ws.Range("L131").Value = 5739
$ws.Range("M131").Value = 300
$ws.Range("N131").Value = -15819

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 13266.333
$ws.Range("I70").Value = 13266.333
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 13266.333
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12996.333
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 13266.333
$ws.Range("I73").Value = 13266.333
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 13266.333
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -12330.333
$ws.Range("N73").ClearContents()
# Row 80
$ws.Range("H80").Value = 2331.889
$ws.Range("I80").Value = 2139.4
$ws.Range("J80").Value = 2572.5
$ws.Range("K80").Value = 2139.4
$ws.Range("L80").Value = 2572.5
$ws.Range("M80").Value = -1141.4
$ws.Range("N80").Value = -4568.5
# Row 83
$ws.Range("H83").Value = 2331.889
$ws.Range("I83").Value = 2139.4
$ws.Range("J83").Value = 2572.5
$ws.Range("K83").Value = 10697
$ws.Range("L83").Value = 12862.5
$ws.Range("M83").Value = -5705
$ws.Range("N83").Value = -22846.5
# Row 126
$ws.Range("H126").Value = 4423.846
$ws.Range("I126").Value = 4864.636
$ws.Range("J126").Value = 1999.5
$ws.Range("K126").Value = 14593.908
$ws.Range("L126").Value = 5998.5
$ws.Range("M126").Value = -12123.908
$ws.Range("N126").Value = -10938.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3404.1428
$ws.Range("I22").Value = 2665.8
$ws.Range("J22").Value = 5250
$ws.Range("K22").Value = 2665.8
$ws.Range("L22").Value = 5250
$ws.Range("M22").Value = -2370.8
$ws.Range("N22").Value = -5840
# Row 27
$ws.Range("H27").Value = 3404.1428
$ws.Range("I27").Value = 2665.8
$ws.Range("J27").Value = 5250
$ws.Range("K27").Value = 2665.8
$ws.Range("L27").Value = 5250
$ws.Range("M27").Value = -2558.8
$ws.Range("N27").Value = -5464
# Row 40
$ws.Range("H40").Value = 1965.7142
$ws.Range("I40").Value = 1867.3158
$ws.Range("K40").Value = 1867.3158
$ws.Range("M40").Value = -1731.3158
# Row 46
$ws.Range("H46").Value = 1550
$ws.Range("I46").Value = 1550
$ws.Range("K46").Value = 1550
$ws.Range("M46").Value = -1362
# Row 93
$ws.Range("H93").Value = 1408.3334
$ws.Range("I93").Value = 1479.4375
$ws.Range("K93").Value = 1479.4375
$ws.Range("M93").Value = -231.4375

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 81
$ws.Range("H81").Value = 2247.6428
$ws.Range("I81").Value = 2451.4167
$ws.Range("J81").Value = 1025
$ws.Range("K81").Value = 4902.8334
$ws.Range("L81").Value = 2050
$ws.Range("M81").Value = -3841.8334
$ws.Range("N81").Value = -4172
# Row 84
$ws.Range("H84").Value = 2247.6428
$ws.Range("I84").Value = 2451.4167
$ws.Range("J84").Value = 1025
$ws.Range("K84").Value = 24514.167
$ws.Range("L84").Value = 10250
$ws.Range("M84").Value = -19210.167
$ws.Range("N84").Value = -20858
